$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# C4 loses its "NA" value (was an inline string "NA", becomes an empty inline string)
$ws.Range("C4").Value = ""

# New row 5 with the data produced by the latest script run.
# Leading apostrophe forces the date-looking text to stay plain text,
# matching how the other date cells (A2:A4) are stored.
$ws.Range("A5").Value = "'2025-03-05"
$ws.Range("B5").Value = "Rien ne nous concerne aujourd'hui !"
$ws.Range("C5").Value = "NA"
$ws.Range("D5").Value = 1
